# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet,
#   with the same column layout as the other quarterly sheets, holding the
#   single fund holding reported for that quarter.
# - Insert a new leading data row in "总计" summarising the new quarter and
#   renumber the running index column for the rows that shift down.

$wb = $excel.ActiveWorkbook

$templateWs = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet, positioned immediately before "总计"
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($beforeSheet)
$q1.Name = "2022-Q1"

# NOTE: worksheet handles are positional - once a sheet is inserted, any
# handle obtained *before* the insert that pointed at (or after) the
# insertion point now resolves to a different sheet. Re-resolve "总计" by
# name after the Add() so later edits land on the right sheet.
$totalWs = $wb.Worksheets.Item("总计")

# Match the sheet-level cosmetics used by the other quarterly sheets.
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Pull the header-row (B1:H1) and the "A" index-column styling from the
# existing 2021-Q4 sheet so the new sheet's look matches its siblings.
$templateWs.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$templateWs.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0

$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "001672"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "国寿安保智慧生活股票"

$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "3.56"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "85.91"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "2.87"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.1022"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row at the top of "总计", pushing the
#    existing rows down and renumbering the running index column.
# ---------------------------------------------------------------------
$totalWs.Rows.Item(2).Insert()
$totalWs.Range("B2:D2").ClearFormats()

$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)
$totalWs.Range("A2").Value = 0

$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 1
$totalWs.Range("D2").Value = 0.1

$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
$totalWs.Range("A7").Value = 5

# Restore the originally active sheet/tab (unchanged by this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
